$d = $word.ActiveDocument

# " Sediment JF743552 Marine sediments " -> " sediment JF743552 marine sediments "
# (lower-case the two capitalized words "Sediment" and "Marine")
$d.Content.Find.Execute("Sediment JF743552 Marine", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "sediment JF743552 marine", 2)

# "Endosymbionts" -> "endosymbionts"
$d.Content.Find.Execute("Endosymbionts", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "endosymbionts", 2)
